# Actualización automática 2025-07-16 14:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
# Row 18 corresponds to ALMEIDA CUATIN JHONATHANN CARLOS / MANCHENO PINO HERVIN SANTIAGO
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M18").Value = 1133.39              # PORCELANATO
$wsGrupo.Range("O18").Value = 74.65000000000001    # SAL SOLUBLE

# --- Sheet 2: VENTA MENSUAL ---
# Row 18 corresponds to the same advisor/client; F = julio (July)
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F18").Value = 1499.94

# Row 32 is the TOTAL row (sum of each column across all advisors/clients)
$wsMensual.Range("F32").Value = 3977.03

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
# Rows for ALMEIDA CUATIN JHONATHANN CARLOS: 15 = PORCELANATO, 17 = SAL SOLUBLE, 18 = TOTAL
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 15 - PORCELANATO
$wsCumpl.Range("D15").Value = 1274.2
$wsCumpl.Range("E15").Value = 22184.62
$wsCumpl.Range("F15").Value = 0.05431645752002872

# Row 17 - SAL SOLUBLE
$wsCumpl.Range("D17").Value = 72.88
$wsCumpl.Range("E17").Value = 1527.12
$wsCumpl.Range("F17").Value = 0.04555

# Row 18 - TOTAL (sum of rows 2-17 for this advisor's groups)
$wsCumpl.Range("D18").Value = 3966.75
$wsCumpl.Range("E18").Value = 29967.96607548726
$wsCumpl.Range("F18").Value = 0.1168935667879473
